# Modification du modèle et du règlement
#
# Target paragraph (style "Explicationdonne") currently reads:
#   "(" + "Film, BD, roman, etc." + ")"
# It must become three runs:
#   "(long-métrage, court-métrage" + ", roman, etc." + ")"
# and a new paragraph (same style) must be inserted right after it:
#   "(cf. les supports acceptés dans le règlement)"

$d = $word.ActiveDocument

# Locate the explanatory paragraph that currently documents the
# "Support" field's accepted formats.
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Film, BD, roman, etc.*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the '(Film, BD, roman, etc.)' paragraph"
}

$pRange = $target.Range
$start = $pRange.Start
$end = $pRange.End

# Replace the paragraph's textual content (everything before the
# paragraph mark) with three runs matching the new wording, built via
# raw WordprocessingML so the run boundaries are explicit.
$body = $d.Range($start, $end - 1)
$newRunsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
    + '<w:r><w:rPr/><w:t>(long-métrage, court-métrage</w:t></w:r>' `
    + '<w:r><w:rPr/><w:t>, roman, etc.</w:t></w:r>' `
    + '<w:r><w:rPr/><w:t>)</w:t></w:r>' `
    + '</w:p>'
[void]$body.InsertXML($newRunsXml)

# Re-fetch the (possibly reseated) paragraph and append a new paragraph
# right after it, with the same "Explicationdonne" style, carrying the
# reference to the rules document.
$target = $d.Paragraphs.Item($target.Index)
$newParaRange = $target.Range.InsertParagraphAfter()

$refreshed = $d.Paragraphs.Item($target.Index + 1)
$refreshed.Range.Text = "(cf. les supports acceptés dans le règlement)"

Write-Output "done"
